# Weekly price-list update: two new price observations are inserted into the
# "Bruselas (repollito)" sheet (one around row 64, one near the end of the
# table), pushing the existing rows below them down by one each time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new record at row 64 (rows 64..93 shift down to 65..94) ---
$ws.Rows.Item(64).Insert()

$ws.Cells.Item(64, 1).Value2  = 6
$ws.Cells.Item(64, 2).Value   = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(64, 3).Value   = "Metropolitana"
$ws.Cells.Item(64, 4).Value2  = 45120
$ws.Cells.Item(64, 5).Value2  = 13
$ws.Cells.Item(64, 6).Value2  = 100112035
$ws.Cells.Item(64, 7).Value   = "Bruselas (repollito)"
$ws.Cells.Item(64, 8).Value   = "Sin especificar"
$ws.Cells.Item(64, 9).Value   = "Primera"
$ws.Cells.Item(64, 10).Value2 = 260
$ws.Cells.Item(64, 11).Value2 = 17000
$ws.Cells.Item(64, 12).Value2 = 18000
$ws.Cells.Item(64, 13).Value2 = 17538
$ws.Cells.Item(64, 14).Value  = "`$/malla 15 kilos"
$ws.Cells.Item(64, 15).Value  = "Provincia de Quillota"
$ws.Cells.Item(64, 16).Value2 = 1169
$ws.Cells.Item(64, 17).Value2 = 15
$ws.Cells.Item(64, 18).Value  = "Hortaliza"

$ws.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Insert second new record at row 93 (current rows 93..94 shift down to 94..95) ---
$ws.Rows.Item(93).Insert()

$ws.Cells.Item(93, 1).Value2  = 6
$ws.Cells.Item(93, 2).Value   = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(93, 3).Value   = "Metropolitana"
$ws.Cells.Item(93, 4).Value2  = 45121
$ws.Cells.Item(93, 5).Value2  = 13
$ws.Cells.Item(93, 6).Value2  = 100112035
$ws.Cells.Item(93, 7).Value   = "Bruselas (repollito)"
$ws.Cells.Item(93, 8).Value   = "Sin especificar"
$ws.Cells.Item(93, 9).Value   = "Primera"
$ws.Cells.Item(93, 10).Value2 = 380
$ws.Cells.Item(93, 11).Value2 = 18000
$ws.Cells.Item(93, 12).Value2 = 20000
$ws.Cells.Item(93, 13).Value2 = 19053
$ws.Cells.Item(93, 14).Value  = "`$/malla 15 kilos"
$ws.Cells.Item(93, 15).Value  = "Provincia de Quillota"
$ws.Cells.Item(93, 16).Value2 = 1270
$ws.Cells.Item(93, 17).Value2 = 15
$ws.Cells.Item(93, 18).Value  = "Hortaliza"

$ws.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
